$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("CS206")

# Remove the header row from CS206; data rows shift up by one (old row 2 -> new
# row 1, ... old row 37 -> new row 36), matching the commit's removal of the
# "Registration Number / Name / Sessional Marks / End Semester Marks /
# Enrollment Status" header and of the now-unused Sessional/End-Semester-Marks
# columns (they were always empty for data rows).
[void]$ws1.Rows(1).Delete()

# Add the new "CS204" sheet right after "CS206" ("Add Loading for Students").
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "CS204"

# Registration number / student-name pairs for the new CS204 roster. A few
# rows reuse existing CS206 students (cs122011, cs122030, cs122043,
# cs131013); the rest are the placeholder "Student#" loading rows.
$colA = @("cs666666", "cs666661", "cs666662", "cs666663", "cs122011", "cs122030", "cs122043", "cs666664", "cs666665", "cs131013")
$colB = @("Student6", "Student1", "Student2", "Student3", "Hafiz Muhammad Haris", "Muhammad Touseef Khan", "Uroosa Shahid", "Student4", "Student5", "Ghulam Rasool")

for ($i = 0; $i -lt $colA.Length; $i++) {
    $ws2.Cells.Item($i + 1, 1).Value = $colA[$i]
}
for ($i = 0; $i -lt $colB.Length; $i++) {
    $ws2.Cells.Item($i + 1, 2).Value = $colB[$i]
}

$ws2.Columns(2).ColumnWidth = 26.14
[void]$ws2.Range("E6").Select()

# Re-activate CS206 so it stays the visible/selected tab, with its reported
# selection, after CS204 is added.
[void]$ws1.Select()
[void]$ws1.Range("C11").Select()
